# raven.docx edit: "added main and modifed"
#
# 1. First paragraph: append two trailing spaces to the existing text run,
#    then append a red-colored annotation "(This is a change - Version for
#    main branch)" as three separate runs (matching the target XML's run
#    split).
# 2. Remove the trailing duplicated/leftover paragraph
#    "ank God almighty, we are free at last."
# 3. Drop a handful of now-unused styles (Heading2/Heading4 and their
#    linked/related styles, plus some web-import leftovers) that the
#    target document no longer carries.

$d = $word.ActiveDocument

# --- 1. Paragraph 1: trailing spaces + red annotation -----------------
$p1 = $d.Paragraphs(1)
$p1Range = $p1.Range
$textEnd = $p1Range.End - 1   # exclude the paragraph mark
$spacer = $d.Range($p1Range.Start, $textEnd)
$spacer.InsertAfter("  ")

$enDash = [char]0x2013

$seg1Start = $d.Paragraphs(1).Range.End - 1
$seg1 = $d.Range($seg1Start, $seg1Start)
$seg1.InsertAfter("(This is a change " + $enDash + " Ve")
$seg1End = $d.Paragraphs(1).Range.End - 1
$d.Range($seg1Start, $seg1End).Font.Color = 255

$seg2Start = $d.Paragraphs(1).Range.End - 1
$seg2 = $d.Range($seg2Start, $seg2Start)
$seg2.InsertAfter("rsion for main branch")
$seg2End = $d.Paragraphs(1).Range.End - 1
$d.Range($seg2Start, $seg2End).Font.Color = 255

$seg3Start = $d.Paragraphs(1).Range.End - 1
$seg3 = $d.Range($seg3Start, $seg3Start)
$seg3.InsertAfter(")")
$seg3End = $d.Paragraphs(1).Range.End - 1
$d.Range($seg3Start, $seg3End).Font.Color = 255

# --- 2. Remove the trailing leftover paragraph -------------------------
$lastIndex = $d.Paragraphs.Count
$d.Paragraphs($lastIndex).Range.Delete()

# --- 3. Remove now-unused styles ---------------------------------------
# Deleting by walking the Styles collection backwards (highest index
# first) so each Delete() call only ever touches the current tail of the
# collection.
$staleStyleNames = @(
    "Heading 2",
    "Heading 4",
    "apple-converted-space",
    "Hyperlink",
    "Heading 2 Char",
    "Heading 4 Char",
    "audio-tool",
    "subscribe",
    "subscribe-more-info",
    "generic-title",
    "podcast-tools__subscribe-links"
)
for ($i = $d.Styles.Count; $i -ge 1; $i--) {
    $style = $d.Styles.Item($i)
    if ($staleStyleNames -contains $style.NameLocal) {
        $style.Delete()
    }
}

Write-Output "done"
